$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1341199514547732
$ws.Range("E2").Value = 10.06518312636318
$ws.Range("F2").Value = 28.68261455400278
